# LOT2040.xlsx edit:
# 1. Remove the row that only held the "8711290 - Elisson Antônio da Costa
#    Romanel" docente value (old row 13, under "Docentes responsáveis:");
#    everything below it shifts up one row.
# 2. Re-point several B/C cells (which keep the same style/role as before
#    but now show different data) to their new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-redundant row (old row 13: blank A / docente name in B+C).
$ws.Rows.Item(13).Delete()

# "Objetivos:" row now shows the docente identification text.
$ws.Cells.Item(10, 2).Value = "8711290 - Elisson Antônio da Costa Romanel"
$ws.Cells.Item(10, 3).Value = "8711290 - Elisson Antônio da Costa Romanel"

# "Programa resumido:" row (new row 13, after the shift) now shows "Semestral".
$ws.Cells.Item(13, 2).Value = "Semestral"
$ws.Cells.Item(13, 3).Value = "Semestral"

# "Programa:" row (new row 15) now shows the activation date (reuse the
# existing "Ativação:" row's text cell via copy/paste-values so it stays a
# plain text shared string instead of Excel auto-converting the "01/01/2019"
# literal into a date-formatted numeric cell).
$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4163)
$ws.Cells.Item(8, 3).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4163)
$excel.CutCopyMode = $false

# "Método:" row (new row 18) now shows the docente identification text.
$ws.Cells.Item(18, 2).Value = "8711290 - Elisson Antônio da Costa Romanel"
$ws.Cells.Item(18, 3).Value = "8711290 - Elisson Antônio da Costa Romanel"

# "Critério:" row (new row 19) now shows the grading-notes text.
$ws.Cells.Item(19, 2).Value = "Notas - N distribuído no semestre. A composição das ""N"" fica critério do docente."
$ws.Cells.Item(19, 3).Value = "Notas - N distribuído no semestre. A composição das ""N"" fica critério do docente."

# "Norma de recuperação:" row (new row 20) now shows the MF formula text.
$ws.Cells.Item(20, 2).Value = "MF = MF = (somatório de N)/número de N (adequando o valor de N, quando houver peso distinto para as Ns)"
$ws.Cells.Item(20, 3).Value = "MF = MF = (somatório de N)/número de N (adequando o valor de N, quando houver peso distinto para as Ns)"

# "Bibliografia:" row (new row 21) now shows the NF/recovery formula text.
$ws.Cells.Item(21, 2).Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0."
$ws.Cells.Item(21, 3).Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0."
